$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2025-04-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-04 Friday", 2) | Out-Null

# Update each multiplication problem cell-by-cell (row, col are 1-based)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "43×34="
$t.Cell(1, 2).Range.Text = "28×50="
$t.Cell(1, 3).Range.Text = "81×14="
$t.Cell(1, 4).Range.Text = "74×88="
$t.Cell(1, 5).Range.Text = "31×85="

$t.Cell(5, 1).Range.Text = "29×95="
$t.Cell(5, 2).Range.Text = "15×34="
$t.Cell(5, 3).Range.Text = "76×97="
$t.Cell(5, 4).Range.Text = "55×71="
$t.Cell(5, 5).Range.Text = "71×33="

$t.Cell(10, 1).Range.Text = "97×62="
$t.Cell(10, 2).Range.Text = "12×55="
$t.Cell(10, 3).Range.Text = "52×65="
$t.Cell(10, 4).Range.Text = "11×83="
$t.Cell(10, 5).Range.Text = "80×68="

$t.Cell(15, 1).Range.Text = "72×55="
$t.Cell(15, 2).Range.Text = "52×93="
$t.Cell(15, 3).Range.Text = "53×26="
$t.Cell(15, 4).Range.Text = "70×16="
$t.Cell(15, 5).Range.Text = "11×29="

$t.Cell(20, 1).Range.Text = "76×40="
$t.Cell(20, 2).Range.Text = "29×68="
$t.Cell(20, 3).Range.Text = "33×77="
$t.Cell(20, 4).Range.Text = "80×74="
$t.Cell(20, 5).Range.Text = "66×41="

